$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.928.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.85%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.286.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.74%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.86"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.55%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.56%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.425"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.849.63"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.41%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.84"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.28%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "68.866.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.80%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.50%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.271.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.24%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.85"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "396.93"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.74"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.56%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.517"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.83%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000121"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.190"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.75"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.99"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.73"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.80%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.02"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.46%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.17"

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.66"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.51"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.55%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.832"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.54%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.44"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.79%  "

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.59"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.41%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.60"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.76%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.56"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.72%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.28%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0691"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "347.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.25%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.627.74"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.77"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0283"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.97%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.15"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.32"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.50%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.07%  "
